$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range('E2').Value = '2026-02-23 07:18:42'
$ws.Range('K2').Value = '0.0 MJ/m2'
$ws.Range('N2').Value = '0.6 °C 6:59 TU'
$ws.Range('O2').Value = '3.1 °C'
$ws.Range('E3').Value = '2026-02-23 07:18:44'
$ws.Range('E4').Value = '2026-02-23 07:18:47'
$ws.Range('H4').Value = '89%'
$ws.Range('E5').Value = '2026-02-23 07:18:50'
$ws.Range('H5').Value = '35%'
$ws.Range('E6').Value = '2026-02-23 07:18:53'
$ws.Range('N6').Value = '7.5 °C 6:52 TU'
$ws.Range('O6').Value = '9.1 °C'
$ws.Range('E7').Value = '2026-02-23 07:18:55'
$ws.Range('H7').Value = '70%'
$ws.Range('K7').Value = '0.0 MJ/m2'
$ws.Range('O7').Value = '11.7 °C'
$ws.Range('E8').Value = '2026-02-23 07:18:58'
$ws.Range('H8').Value = '56%'
$ws.Range('E9').Value = '2026-02-23 07:19:01'
$ws.Range('H9').Value = '91%'
$ws.Range('O9').Value = '6.6 °C'
$ws.Range('E10').Value = '2026-02-23 07:19:04'
$ws.Range('H10').Value = '98%'
$ws.Range('I10').Value = '0.1 mm'
$ws.Range('K10').Value = '0.0 MJ/m2'
$ws.Range('N10').Value = '3.1 °C 6:36 TU'
$ws.Range('E11').Value = '2026-02-23 07:19:06'
$ws.Range('N11').Value = '1.3 °C 6:54 TU'
$ws.Range('O11').Value = '2.6 °C'
$ws.Range('E12').Value = '2026-02-23 07:19:09'
$ws.Range('O12').Value = '5.2 °C'
$ws.Range('E13').Value = '2026-02-23 07:19:11'
$ws.Range('H13').Value = '91%'
$ws.Range('J13').Value = '1032.3 hPa'
$ws.Range('N13').Value = '-3.6 °C 6:45 TU'
$ws.Range('O13').Value = '-1.4 °C'
$ws.Range('E14').Value = '2026-02-23 07:19:14'
$ws.Range('E15').Value = '2026-02-23 07:19:17'
$ws.Range('H15').Value = '89%'
$ws.Range('N15').Value = '4.5 °C 6:41 TU'
$ws.Range('O15').Value = '6.4 °C'
$ws.Range('E16').Value = '2026-02-23 07:19:19'
$ws.Range('H16').Value = '17%'
$ws.Range('E17').Value = '2026-02-23 07:19:22'
$ws.Range('N17').Value = '6.3 °C 6:51 TU'
$ws.Range('E18').Value = '2026-02-23 07:19:25'
$ws.Range('N18').Value = '1.1 °C 6:46 TU'
$ws.Range('O18').Value = '2.7 °C'
$ws.Range('E19').Value = '2026-02-23 07:19:28'
$ws.Range('N19').Value = '7.7 °C 6:35 TU'
$ws.Range('O19').Value = '9.6 °C'
$ws.Range('E20').Value = '2026-02-23 07:19:30'
$ws.Range('H20').Value = '35%'
$ws.Range('E21').Value = '2026-02-23 07:19:33'
$ws.Range('J21').Value = '1029.4 hPa'
$ws.Range('L21').Value = '8.3 km/h - 338º 6:45 TU'
$ws.Range('N21').Value = '1.1 °C 6:30 TU'
$ws.Range('O21').Value = '3.2 °C'
$ws.Range('E22').Value = '2026-02-23 07:19:36'
$ws.Range('O22').Value = '1.9 °C'
$ws.Range('E23').Value = '2026-02-23 07:19:38'
$ws.Range('H23').Value = '24%'
$ws.Range('L23').Value = '43.9 km/h - 323º 6:36 TU'
$ws.Range('E24').Value = '2026-02-23 07:19:41'
$ws.Range('N24').Value = '0.1 °C 6:39 TU'
$ws.Range('O24').Value = '2.0 °C'
$ws.Range('E25').Value = '2026-02-23 07:19:44'
$ws.Range('E26').Value = '2026-02-23 07:19:46'
$ws.Range('J26').Value = '1026.4 hPa'
$ws.Range('K26').Value = '0.0 MJ/m2'
$ws.Range('L26').Value = '14.4 km/h - 26º 6:33 TU'
$ws.Range('M26').Value = '7.6 °C 6:51 TU'
$ws.Range('O26').Value = '6.5 °C'
$ws.Range('E27').Value = '2026-02-23 07:19:49'
$ws.Range('E28').Value = '2026-02-23 07:19:52'
$ws.Range('O28').Value = '3.4 °C'
$ws.Range('E29').Value = '2026-02-23 07:19:54'
$ws.Range('K29').Value = '0.0 MJ/m2'
$ws.Range('O29').Value = '4.0 °C'
$ws.Range('E30').Value = '2026-02-23 07:19:56'
$ws.Range('N30').Value = '7.4 °C 6:48 TU'
$ws.Range('E31').Value = '2026-02-23 07:19:59'
$ws.Range('J31').Value = '1024.2 hPa'
$ws.Range('N31').Value = '13.6 °C 6:59 TU'
$ws.Range('E32').Value = '2026-02-23 07:20:02'
$ws.Range('E33').Value = '2026-02-23 07:20:05'
$ws.Range('H33').Value = '66%'
$ws.Range('E34').Value = '2026-02-23 07:20:07'
$ws.Range('H34').Value = '45%'
$ws.Range('M34').Value = '5.4 °C 6:36 TU'
$ws.Range('O34').Value = '2.3 °C'
$ws.Range('E35').Value = '2026-02-23 07:20:10'
$ws.Range('J35').Value = '1026.4 hPa'
$ws.Range('L35').Value = '42.8 km/h - 281º 6:33 TU'
$ws.Range('E36').Value = '2026-02-23 07:20:13'
$ws.Range('E37').Value = '2026-02-23 07:20:15'
$ws.Range('J37').Value = '1029.6 hPa'
$ws.Range('N37').Value = '0.8 °C 6:34 TU'
$ws.Range('O37').Value = '3.3 °C'
$ws.Range('E38').Value = '2026-02-23 07:20:18'
$ws.Range('H38').Value = '74%'
$ws.Range('K38').Value = '0.0 MJ/m2'
$ws.Range('E39').Value = '2026-02-23 07:20:20'
$ws.Range('H39').Value = '23%'
$ws.Range('K39').Value = '0.0 MJ/m2'
$ws.Range('L39').Value = '33.5 km/h - 335º 6:51 TU'
$ws.Range('O39').Value = '3.5 °C'
$ws.Range('E40').Value = '2026-02-23 07:20:23'
$ws.Range('H40').Value = '92%'
$ws.Range('J40').Value = '1029.8 hPa'
$ws.Range('N40').Value = '-0.2 °C 6:31 TU'
$ws.Range('O40').Value = '1.7 °C'
$ws.Range('E41').Value = '2026-02-23 07:20:25'
$ws.Range('K41').Value = '0.0 MJ/m2'
$ws.Range('E42').Value = '2026-02-23 07:20:28'
$ws.Range('E43').Value = '2026-02-23 07:20:30'
$ws.Range('H43').Value = '96%'
$ws.Range('K43').Value = '0.0 MJ/m2'
$ws.Range('O43').Value = '3.5 °C'
$ws.Range('E44').Value = '2026-02-23 07:20:33'
$ws.Range('H44').Value = '36%'
$ws.Range('K44').Value = '0.0 MJ/m2'
$ws.Range('E45').Value = '2026-02-23 07:20:36'
$ws.Range('H45').Value = '65%'
$ws.Range('J45').Value = '1030.4 hPa'
$ws.Range('E46').Value = '2026-02-23 07:20:38'
$ws.Range('N46').Value = '0.3 °C 6:47 TU'
$ws.Range('O46').Value = '1.8 °C'
